$p = $ppt.ActivePresentation
$newDate = "09-Apr-19"

# --- 1. Update the cached "datetimeFigureOut" date placeholders -------------
# These live on the Slide Master, every slide Layout, and the Notes Master.

$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$notesMaster = $p.NotesMaster
foreach ($shp in $notesMaster.Shapes) {
    if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- 2. Rename the AddressBook-era API calls to their Task equivalents ------

$slide = $p.Slides.Item(1)

# "deletePerson(p)" -> "deleteTask(t)"
$deleteShape = $slide.Shapes.Item("TextBox 28")
$tr = $deleteShape.TextFrame.TextRange
$run1 = $tr.Characters(1, 12)
$run1.Text = "deleteTask"
$tr = $deleteShape.TextFrame.TextRange
$run2 = $tr.Characters(11, 3)
$run2.Text = "(t)"

# "saveAddressBook(AddressBook)" -> "saveTaskBook(TaskBook)"
$saveShape = $slide.Shapes.Item("TextBox 73")
$tr = $saveShape.TextFrame.TextRange
$run3 = $tr.Characters(1, 15)
$run3.Text = "saveTaskBook"
$tr = $saveShape.TextFrame.TextRange
$run4 = $tr.Characters(14, 11)
$run4.Text = "TaskBook"
